$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.22594800312889163
$ws.Range("B1").Value = 0.2257775568961975
$ws.Range("A2").Value = -0.20367121235961694
$ws.Range("B2").Value = 0.20299410486823621
$ws.Range("A3").Value = -0.11124527720195942
$ws.Range("B3").Value = 0.11093568000425336
$ws.Range("A4").Value = -0.1029356800225969
$ws.Range("B4").Value = 0.10247341700745238
$ws.Range("A5").Value = -0.099473417018081989
$ws.Range("B5").Value = 0.097902837865131431
$ws.Range("A6").Value = -0.047193694216497661
$ws.Range("B6").Value = 0.046725314151691677
$ws.Range("A7").Value = -0.036725314177425972
$ws.Range("B7").Value = 0.036613820809900943
$ws.Range("A8").Value = -0.026613820836378199
$ws.Range("B8").Value = 0.02642359506033376
$ws.Range("A9").Value = -0.024423595073642002
$ws.Range("B9").Value = 0.024269930043501642
$ws.Range("A10").Value = -0.022269930057625231
$ws.Range("B10").Value = 0.02225976793762463
$ws.Range("A11").Value = -0.019259767953597517
$ws.Range("B11").Value = 0.01924393606816821
$ws.Range("A12").Value = -0.015743936085210297
$ws.Range("B12").Value = 0.015637739309380283
$ws.Range("A13").Value = -0.017172187527815019
$ws.Range("B13").Value = 0.01708323456229266
$ws.Range("A14").Value = -0.0090832345882789767
$ws.Range("B14").Value = 0.0090539229138801858
$ws.Range("A15").Value = -0.0080539229277638569
$ws.Range("B15").Value = 0.0080349716956913397
$ws.Range("A16").Value = -0.0060349717115113499
$ws.Range("B16").Value = 0.006003747258722214
$ws.Range("A17").Value = -0.0040037472748082337
$ws.Range("B17").Value = 0.0039999999803264075
$ws.Range("A18").Value = -0.016106344547644369
$ws.Range("B18").Value = 0.016091941873813909
$ws.Range("A19").Value = -0.012091941881424262
$ws.Range("B19").Value = 0.012017209502860116
$ws.Range("A20").Value = -0.0080172095109052322
$ws.Range("B20").Value = 0.0080057216971578526
$ws.Range("A21").Value = -0.0040057217052904548
$ws.Range("B21").Value = 0.003999999991803449
$ws.Range("A22").Value = -0.087748827676540131
$ws.Range("B22").Value = 0.087170243019917493
$ws.Range("A23").Value = -0.075636629491797613
$ws.Range("B23").Value = 0.074563893705533957
$ws.Range("A24").Value = -0.020099660266230224
$ws.Range("B24").Value = 0.019999999958019821
$ws.Range("A25").Value = -0.048709143660063958
$ws.Range("B25").Value = 0.048674389228363424
$ws.Range("A26").Value = -0.046174389240887237
$ws.Range("B26").Value = 0.0461332876537206
$ws.Range("A27").Value = -0.043633287666605458
$ws.Range("B27").Value = 0.043411705595068018
$ws.Range("A28").Value = -0.041411705608473071
$ws.Range("B28").Value = 0.041274662876591428
$ws.Range("A29").Value = -0.034274662899739461
$ws.Range("B29").Value = 0.03424543224528076
$ws.Range("A30").Value = 0.025754567638212489
$ws.Range("B30").Value = -0.025797868836482607
$ws.Range("A31").Value = 0.032797868813577935
$ws.Range("B31").Value = -0.032826981644157982
$ws.Range("A32").Value = -0.0040016224004357781
$ws.Range("B32").Value = 0.0039999999825788279
$ws.Columns("B").ColumnWidth = 13.8
